$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (row 3) duplicating the existing row 2 match record.
$ws.Range("A3").Value = " Dubai (DSC)"
$ws.Range("B3").Value = " September 24 2020"
$ws.Range("C3").Value = "Kings XI won by 97 runs"
$ws.Range("D3").Value = "Royal Challengers Bangalore"
$ws.Range("E3").Value = "Kings XI Punjab"
$ws.Range("F3").Value = "Dale Steyn "

# These columns look numeric but must stay text (matching row 2's storage),
# so force text format before assigning, then drop back to the default style.
$ws.Range("G3:K3").NumberFormat = "@"
$ws.Range("G3").Value = "1"
$ws.Range("H3").Value = "2"
$ws.Range("I3").Value = "0"
$ws.Range("J3").Value = "0"
$ws.Range("K3").Value = "50.00"
$ws.Range("G3:K3").Style = "Normal"
